# "trying to get parameter excel to be read"
# Add new rows of parameters to the "fixed parameters" sheet and a new
# "datetime" header column to the "dynamic parameters" sheet, then make
# "dynamic parameters" the active tab/selection.

$wb = $excel.ActiveWorkbook

$wsFixed = $wb.Worksheets.Item("fixed parameters")
$wsDynamic = $wb.Worksheets.Item("dynamic parameters")

# --- "fixed parameters" sheet: new parameter rows ---
$wsFixed.Range("A2").Value = "temperature"
$wsFixed.Range("B2").Value = 20

$wsFixed.Range("A3").Value = "var A"
$wsFixed.Range("B3").Value = "A"

$wsFixed.Range("A4").Value = "var B"
$wsFixed.Range("B4").Value = 33

$wsFixed.Range("A5").Value = "var C"

# --- "dynamic parameters" sheet: new header column ---
$wsDynamic.Range("A1").Value = "datetime"

# --- selections ---
$wsFixed.Range("B5").Select() | Out-Null
$wsDynamic.Range("A2").Select() | Out-Null

# "dynamic parameters" becomes the active/visible tab
$wsDynamic.Activate()
